# Generate Report for handback
# Row 2 on each language sheet holds the most recent handback record
# (937e6c38-... file); re-running the report generation refreshes its
# "Correspond Handoff Datetime" (col D) and "Correspond Handback DateTime"
# (col G) with the newest timestamps, while row 3 (the f8972293-... file)
# keeps its original values.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-26 11:53:19"
$wsZhCn.Range("G2").Value = "2016-01-26 11:54:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-26 11:53:30"
$wsDeDe.Range("G2").Value = "2016-01-26 11:54:28"
